$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.085.09"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "2.240.50"
$ws.Range("E3").Value = "  -4.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.95%  "

$ws.Range("E7").Value = "  -2.01%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0798"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.57%  "

$ws.Range("E13").Value = "  -2.04%  "

$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").Value = "2.591.25"
$ws.Range("E15").Value = "  -3.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.73%  "

$ws.Range("D17").Value = "2.235.00"
$ws.Range("E17").Value = "  -4.93%  "

$ws.Range("E18").Value = "  -3.14%  "

$ws.Range("D19").Value = "40.019.58"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").Value = "0.0₃0894"

$ws.Range("E21").Value = "  -3.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("E34").Value = "  -3.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0713"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  -3.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.74%  "

$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("E39").Value = "  +0.73%  "

$ws.Range("E40").Value = "  -1.86%  "

$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("D43").Value = "1.958.59"
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("E45").Value = "  +3.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.01%  "

$ws.Range("E48").Value = "  -0.91%  "

$ws.Range("D49").Value = "2.457.26"
$ws.Range("E49").Value = "  -4.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.93%  "

$ws.Range("E51").Value = "  +8.41%  "
